# Add a 4th and 5th bridge domain (BD) and their corresponding EPGs,
# normalize the BD description text, reorder the EPG sheet ahead of
# SUBNET, and move the active tab to LACP_POLICY.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. BD sheet: normalize existing descriptions & append two new rows
# ---------------------------------------------------------------------
$bd = $wb.Worksheets.Item("BD")

# New BD name cells first (creates the new shared strings in the same
# order the original author typed them).
$bd.Cells.Item(5, 2).Value = "mark_fourth_bd_for_subnet"
$bd.Cells.Item(6, 2).Value = "mark_fifth_bd_for_subnet"

# Descriptions: 3rd row gets cleaned up (drop stray trailing "1"), then
# the two brand new rows, then the 2nd and 1st row descriptions are
# normalized to match the "Nth bridge domain" wording.
$bd.Cells.Item(4, 3).Value = "This 3rd  bridge domain is created by the Terraform ACI provider"
$bd.Cells.Item(5, 3).Value = "This 4th  bridge domain is created by the Terraform ACI provider"
$bd.Cells.Item(6, 3).Value = "This 5th  bridge domain is created by the Terraform ACI provider"
$bd.Cells.Item(3, 3).Value = "This 2nd  bridge domain is created by the Terraform ACI provider"
$bd.Cells.Item(2, 3).Value = "This 1st bridge domain is created by the Terraform ACI provider"

# type column for the new rows
$bd.Cells.Item(5, 1).Value = "bd"
$bd.Cells.Item(6, 1).Value = "bd"

$bd.Range("C26").Select()

# ---------------------------------------------------------------------
# 2. Reorder sheets: EPG moves ahead of SUBNET
# ---------------------------------------------------------------------
$epgBeforeMove = $wb.Worksheets.Item("EPG")
$subnet = $wb.Worksheets.Item("SUBNET")
$epgBeforeMove.Move($subnet)

# Re-fetch the EPG worksheet by name: after Move() the old handle now
# tracks whichever sheet ended up in that tab slot (SUBNET), not EPG.
$epg = $wb.Worksheets.Item("EPG")

# ---------------------------------------------------------------------
# 3. EPG sheet: append the two new EPGs for the new BDs
# ---------------------------------------------------------------------
$epg.Cells.Item(5, 1).Value = "epg"
$epg.Cells.Item(5, 2).Value = "mark_fouth_epg"
$epg.Cells.Item(5, 3).Value = "mark_fourth_bd_for_subnet"
$epg.Cells.Item(5, 4).Value = "prod_app_profile"

$epg.Cells.Item(6, 1).Value = "epg"
$epg.Cells.Item(6, 2).Value = "mark_fifth_epg"
$epg.Cells.Item(6, 3).Value = "mark_fifth_bd_for_subnet"
$epg.Cells.Item(6, 4).Value = "prod_app_profile"

$epg.Range("B6").Select()

# ---------------------------------------------------------------------
# 4. Move the active tab from TENANT to LACP_POLICY
# ---------------------------------------------------------------------
$lacp = $wb.Worksheets.Item("LACP_POLICY")
$lacp.Activate()
